$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new test-case row (row 8) ---
# Copy the formatting of the last existing data row (row 7) down to row 8
# so the new row picks up the same borders/style as the rest of the table.
$ws.Range("A7:F7").Copy()
$ws.Range("A8:F8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the new row's values. Column order below intentionally sets
# B8 before A8 so the shared-string table grows in the same order as the
# target workbook (B8's string becomes index 21, A8's becomes index 22).
$ws.Range("B8").Value = "Add New Event in Calender"
$ws.Range("A8").Value = "EXL_CorporateLensHomePage_CalendarEvent"
$ws.Range("C8").Value = "N"
$ws.Range("D8").Value = "Y"
$ws.Range("F8").Value = "Sprint1"

# --- Extend the data validation ranges to cover the new row ---
$ws.Range("C2:D7").Validation.Delete() | Out-Null
$ws.Range("F2:F7").Validation.Delete() | Out-Null
$ws.Range("C2:D8").Validation.Add(3, 1, 1, """Y,N""") | Out-Null
$ws.Range("F2:F8").Validation.Add(3, 1, 1, """Sprint1,Sprint2,Sprint3,Sprint4,Sprint5,Sprint6,Sprint7,Sprint8,Sprint9,Sprint10""") | Out-Null

# --- Tidy up the view state ---
# Move the selection back to the top of the sheet (the author's saved file
# no longer shows a stale A7 selection sitting on the old last row).
$ws.Range("A1").Select() | Out-Null
